$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 119 (shifts existing rows 119-131 down to 120-132)
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new weekly data record.
# Columns A,B,C,E,F,G,H,I,J,K,L,R share the same values as the surrounding rows.
$ws.Cells.Item(119, 1).Value = 10
$ws.Cells.Item(119, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(119, 3).Value = "La Araucanía"
$ws.Cells.Item(119, 4).Value = 44449
$ws.Cells.Item(119, 5).Value = 9
$ws.Cells.Item(119, 6).Value = "Fruta"
$ws.Cells.Item(119, 7).Value = 100102
$ws.Cells.Item(119, 8).Value = "Cítricos"
$ws.Cells.Item(119, 9).Value = 100102006
$ws.Cells.Item(119, 10).Value = "Pomelo"
$ws.Cells.Item(119, 11).Value = "Start Ruby"
$ws.Cells.Item(119, 12).Value = "Primera"
$ws.Cells.Item(119, 13).Value = 50
$ws.Cells.Item(119, 14).Value = 12000
$ws.Cells.Item(119, 15).Value = 12000
$ws.Cells.Item(119, 16).Value = 12000
$ws.Cells.Item(119, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(119, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(119, 19).Value = 800
$ws.Cells.Item(119, 20).Value = 15
